$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 3
$ws_ALC.Range("H3").Value = 41440
$ws_ALC.Range("J3").Value = 41440
$ws_ALC.Range("L3").Value = 41440
$ws_ALC.Range("N3").Value = -41668

# ALC row 98
$ws_ALC.Range("H98").Value = 2402.2856
$ws_ALC.Range("I98").Value = 2301.6667
$ws_ALC.Range("K98").Value = 2301.6667
$ws_ALC.Range("M98").Value = -803.6667000000002

# ALC row 101
$ws_ALC.Range("H101").Value = 833.3333
$ws_ALC.Range("I101").Value = 1500
$ws_ALC.Range("J101").Value = 700
$ws_ALC.Range("K101").Value = 4500
$ws_ALC.Range("L101").Value = 2100
$ws_ALC.Range("M101").Value = -2878
$ws_ALC.Range("N101").Value = -5344

# ALC row 102
$ws_ALC.Range("H102").Value = 41440
$ws_ALC.Range("J102").Value = 41440
$ws_ALC.Range("L102").Value = 41440
$ws_ALC.Range("N102").Value = -47930

# ALC row 105
$ws_ALC.Range("H105").Value = 31200
$ws_ALC.Range("J105").Value = 31200
$ws_ALC.Range("L105").Value = 31200
$ws_ALC.Range("N105").Value = -38188

# ALC row 107
$ws_ALC.Range("H107").Value = 12500687
$ws_ALC.Range("I107").Value = 13158355
$ws_ALC.Range("K107").Value = 13158355
$ws_ALC.Range("M107").Value = -13156435

# ALC row 108
$ws_ALC.Range("H108").Value = 43000
$ws_ALC.Range("J108").Value = 43000
$ws_ALC.Range("L108").Value = 43000
$ws_ALC.Range("N108").Value = -50680

# ALC row 110
$ws_ALC.Range("H110").Value = 39351
$ws_ALC.Range("J110").Value = 39351
$ws_ALC.Range("L110").Value = 39351
$ws_ALC.Range("N110").Value = -47531

# ALC row 111
$ws_ALC.Range("H111").Value = 767.8182
$ws_ALC.Range("J111").Value = 866
$ws_ALC.Range("L111").Value = 2598
$ws_ALC.Range("N111").Value = -8732

# ALC row 112
$ws_ALC.Range("H112").Value = 12699507
$ws_ALC.Range("I112").Value = 300
$ws_ALC.Range("J112").Value = 13606593
$ws_ALC.Range("K112").Value = 900
$ws_ALC.Range("L112").Value = 40819779
$ws_ALC.Range("M112").Value = 208
$ws_ALC.Range("N112").Value = -40821995

# ALC row 113
$ws_ALC.Range("H113").Value = 2627.6924
$ws_ALC.Range("I113").Value = 2351.111
$ws_ALC.Range("J113").Value = 3250
$ws_ALC.Range("K113").Value = 2351.111
$ws_ALC.Range("L113").Value = 3250
$ws_ALC.Range("M113").Value = 902.8890000000001
$ws_ALC.Range("N113").Value = -9758

# ALC row 115
$ws_ALC.Range("H115").Value = 7748.2666
$ws_ALC.Range("I115").Value = 11580.444
$ws_ALC.Range("J115").Value = 2000
$ws_ALC.Range("K115").Value = 34741.33199999999
$ws_ALC.Range("L115").Value = 6000
$ws_ALC.Range("M115").Value = -33174.33199999999
$ws_ALC.Range("N115").Value = -9134

# ALC row 116
$ws_ALC.Range("H116").Value = 24002.2
$ws_ALC.Range("I116").Value = 28751.25
$ws_ALC.Range("J116").Value = 5006
$ws_ALC.Range("K116").Value = 28751.25
$ws_ALC.Range("L116").Value = 5006
$ws_ALC.Range("M116").Value = -25309.25
$ws_ALC.Range("N116").Value = -11890

# ALC row 117
$ws_ALC.Range("H117").Value = 0
$ws_ALC.Range("J117").Value = 0
$ws_ALC.Range("L117").Value = 0
$ws_ALC.Range("N117").Value = $null

# ALC row 118
$ws_ALC.Range("H118").Value = 757.64703
$ws_ALC.Range("I118").Value = 631.4
$ws_ALC.Range("J118").Value = 1704.5
$ws_ALC.Range("K118").Value = 1894.2
$ws_ALC.Range("L118").Value = 5113.5
$ws_ALC.Range("M118").Value = -237.1999999999998
$ws_ALC.Range("N118").Value = -8427.5

# ALC row 122
$ws_ALC.Range("H122").Value = 2402.2856
$ws_ALC.Range("I122").Value = 2301.6667
$ws_ALC.Range("K122").Value = 6905.000100000001
$ws_ALC.Range("M122").Value = -4455.000100000001

# ALC row 125
$ws_ALC.Range("H125").Value = 5060
$ws_ALC.Range("J125").Value = 768
$ws_ALC.Range("L125").Value = 6912
$ws_ALC.Range("N125").Value = -11832

# ALC row 128
$ws_ALC.Range("H128").Value = 76778.336
$ws_ALC.Range("J128").Value = 76778.336
$ws_ALC.Range("L128").Value = 76778.336
$ws_ALC.Range("N128").Value = -86738.336

# ALC row 129
$ws_ALC.Range("H129").Value = 1216.6786
$ws_ALC.Range("I129").Value = 682
$ws_ALC.Range("K129").Value = 2046
$ws_ALC.Range("M129").Value = 2954

# ALC row 131
$ws_ALC.Range("H131").Value = 5186.875
$ws_ALC.Range("I131").Value = 1498.3334
$ws_ALC.Range("J131").Value = 7400
$ws_ALC.Range("K131").Value = 4495.0002
$ws_ALC.Range("L131").Value = 22200
$ws_ALC.Range("M131").Value = 544.9997999999996
$ws_ALC.Range("N131").Value = -32280

# ARM row 8
$ws_ARM.Range("H8").Value = 10000
$ws_ARM.Range("I8").Value = 0
$ws_ARM.Range("K8").Value = 0
$ws_ARM.Range("M8").Value = $null

# ARM row 92
$ws_ARM.Range("H92").Value = 28200
$ws_ARM.Range("J92").Value = 28200
$ws_ARM.Range("L92").Value = 28200
$ws_ARM.Range("N92").Value = -33192

# ARM row 94
$ws_ARM.Range("H94").Value = 33815
$ws_ARM.Range("J94").Value = 33815
$ws_ARM.Range("L94").Value = 33815
$ws_ARM.Range("N94").Value = -35617

# BSM row 22
$ws_BSM.Range("H22").Value = 427
$ws_BSM.Range("I22").Value = 383.33334
$ws_BSM.Range("J22").Value = 689
$ws_BSM.Range("K22").Value = 383.33334
$ws_BSM.Range("L22").Value = 689
$ws_BSM.Range("M22").Value = -210.33334
$ws_BSM.Range("N22").Value = -1035

# CRP row 88
$ws_CRP.Range("H88").Value = 36299.855
$ws_CRP.Range("J88").Value = 36299.855
$ws_CRP.Range("L88").Value = 36299.855
$ws_CRP.Range("N88").Value = -37111.855

# CRP row 91
$ws_CRP.Range("H91").Value = 36299.855
$ws_CRP.Range("J91").Value = 36299.855
$ws_CRP.Range("L91").Value = 36299.855
$ws_CRP.Range("N91").Value = -39107.855

# CUL row 56
$ws_CUL.Range("H56").Value = 7500
$ws_CUL.Range("I56").Value = 7500
$ws_CUL.Range("K56").Value = 7500
$ws_CUL.Range("M56").Value = -6970

# CUL row 130
$ws_CUL.Range("H130").Value = 6757.846
$ws_CUL.Range("I130").Value = 2824.5
$ws_CUL.Range("J130").Value = 7207.3716
$ws_CUL.Range("K130").Value = 8473.5
$ws_CUL.Range("L130").Value = 21622.1148
$ws_CUL.Range("M130").Value = -3453.5
$ws_CUL.Range("N130").Value = -31662.1148

# GSM row 122
$ws_GSM.Range("H122").Value = 62640460
$ws_GSM.Range("I122").Value = 133104310
$ws_GSM.Range("J122").Value = 5921.6665
$ws_GSM.Range("K122").Value = 399312930
$ws_GSM.Range("L122").Value = 17764.9995
$ws_GSM.Range("M122").Value = -399310480
$ws_GSM.Range("N122").Value = -22664.9995

# GSM row 133
$ws_GSM.Range("H133").Value = 41091.25
$ws_GSM.Range("J133").Value = 41091.25
$ws_GSM.Range("L133").Value = 41091.25
$ws_GSM.Range("N133").Value = -51211.25

# LTW row 7
$ws_LTW.Range("H7").Value = 2354.8462
$ws_LTW.Range("I7").Value = 1780.8
$ws_LTW.Range("K7").Value = 1780.8
$ws_LTW.Range("M7").Value = -1668.8

# LTW row 22
$ws_LTW.Range("H22").Value = 2249.5
$ws_LTW.Range("J22").Value = 2249.5
$ws_LTW.Range("L22").Value = 2249.5
$ws_LTW.Range("N22").Value = -2839.5

# LTW row 27
$ws_LTW.Range("H27").Value = 2249.5
$ws_LTW.Range("J27").Value = 2249.5
$ws_LTW.Range("L27").Value = 2249.5
$ws_LTW.Range("N27").Value = -2463.5

# LTW row 40
$ws_LTW.Range("H40").Value = 66669588
$ws_LTW.Range("I40").Value = 71431544
$ws_LTW.Range("J40").Value = 2200
$ws_LTW.Range("K40").Value = 71431544
$ws_LTW.Range("L40").Value = 2200
$ws_LTW.Range("M40").Value = -71431408
$ws_LTW.Range("N40").Value = -2472

# LTW row 61
$ws_LTW.Range("H61").Value = 1926.0769
$ws_LTW.Range("J61").Value = 2252.5
$ws_LTW.Range("L61").Value = 2252.5
$ws_LTW.Range("N61").Value = -2656.5

# LTW row 113
$ws_LTW.Range("H113").Value = 1926.0769
$ws_LTW.Range("J113").Value = 2252.5
$ws_LTW.Range("L113").Value = 2252.5
$ws_LTW.Range("N113").Value = -6592.5

# LTW row 122
$ws_LTW.Range("H122").Value = 4796857
$ws_LTW.Range("I122").Value = 10215939
$ws_LTW.Range("J122").Value = 1003500
$ws_LTW.Range("K122").Value = 30647817
$ws_LTW.Range("L122").Value = 3010500
$ws_LTW.Range("M122").Value = -30645367
$ws_LTW.Range("N122").Value = -3015400

# LTW row 126
$ws_LTW.Range("H126").Value = 2354.8462
$ws_LTW.Range("I126").Value = 1780.8
$ws_LTW.Range("K126").Value = 5342.4
$ws_LTW.Range("M126").Value = -2872.4

# WVR row 122
$ws_WVR.Range("H122").Value = 1550.409
$ws_WVR.Range("I122").Value = 1538.2778
$ws_WVR.Range("J122").Value = 1605
$ws_WVR.Range("K122").Value = 4614.8334
$ws_WVR.Range("L122").Value = 4815
$ws_WVR.Range("M122").Value = -2164.8334
$ws_WVR.Range("N122").Value = -9715
